$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 64
$ws.Cells.Item(64, 8).Value = 32823.676
$ws.Cells.Item(64, 9).Value = 3140.4
$ws.Cells.Item(64, 10).Value = 56257.844
$ws.Cells.Item(64, 11).Value = 3140.4
$ws.Cells.Item(64, 12).Value = 56257.844
$ws.Cells.Item(64, 13).Value = -2892.4
$ws.Cells.Item(64, 14).Value = -56753.844

# ALC row 67
$ws.Cells.Item(67, 8).Value = 32823.676
$ws.Cells.Item(67, 9).Value = 3140.4
$ws.Cells.Item(67, 10).Value = 56257.844
$ws.Cells.Item(67, 11).Value = 3140.4
$ws.Cells.Item(67, 12).Value = 56257.844
$ws.Cells.Item(67, 13).Value = -2282.4
$ws.Cells.Item(67, 14).Value = -57973.844

# ALC row 74
$ws.Cells.Item(74, 8).Value = 3437.08
$ws.Cells.Item(74, 9).Value = 3316.1875
$ws.Cells.Item(74, 10).Value = 3652
$ws.Cells.Item(74, 11).Value = 3316.1875
$ws.Cells.Item(74, 12).Value = 3652
$ws.Cells.Item(74, 13).Value = -2380.1875
$ws.Cells.Item(74, 14).Value = -5524

# ALC row 77
$ws.Cells.Item(77, 8).Value = 3437.08
$ws.Cells.Item(77, 9).Value = 3316.1875
$ws.Cells.Item(77, 10).Value = 3652
$ws.Cells.Item(77, 11).Value = 16580.9375
$ws.Cells.Item(77, 12).Value = 18260
$ws.Cells.Item(77, 13).Value = -11900.9375
$ws.Cells.Item(77, 14).Value = -27620

# ALC row 129
$ws.Cells.Item(129, 8).Value = 1001.40814
$ws.Cells.Item(129, 9).Value = 346.66666
$ws.Cells.Item(129, 10).Value = 1044.1086
$ws.Cells.Item(129, 11).Value = 1039.99998
$ws.Cells.Item(129, 12).Value = 3132.3258
$ws.Cells.Item(129, 13).Value = 3960.00002
$ws.Cells.Item(129, 14).Value = -13132.3258

# ALC row 138
$ws.Cells.Item(138, 8).Value = 2683.7017
$ws.Cells.Item(138, 9).Value = 1126.0769
$ws.Cells.Item(138, 10).Value = 3990.0967
$ws.Cells.Item(138, 11).Value = 3378.2307
$ws.Cells.Item(138, 12).Value = 11970.2901
$ws.Cells.Item(138, 13).Value = 1761.7693
$ws.Cells.Item(138, 14).Value = -22250.2901

$ws = $wb.Worksheets.Item("ARM")
# ARM row 74
$ws.Cells.Item(74, 8).Value = 6635.1665
$ws.Cells.Item(74, 9).Value = 11073.1
$ws.Cells.Item(74, 10).Value = 1087.75
$ws.Cells.Item(74, 11).Value = 11073.1
$ws.Cells.Item(74, 12).Value = 1087.75
$ws.Cells.Item(74, 13).Value = -10199.1
$ws.Cells.Item(74, 14).Value = -2835.75

# ARM row 77
$ws.Cells.Item(77, 8).Value = 6635.1665
$ws.Cells.Item(77, 9).Value = 11073.1
$ws.Cells.Item(77, 10).Value = 1087.75
$ws.Cells.Item(77, 11).Value = 55365.5
$ws.Cells.Item(77, 12).Value = 5438.75
$ws.Cells.Item(77, 13).Value = -50997.5
$ws.Cells.Item(77, 14).Value = -14174.75

$ws = $wb.Worksheets.Item("BSM")
# BSM row 86
$ws.Cells.Item(86, 8).Value = 2129.3096
$ws.Cells.Item(86, 9).Value = 2062.8108
$ws.Cells.Item(86, 10).Value = 2621.4
$ws.Cells.Item(86, 11).Value = 2062.8108
$ws.Cells.Item(86, 12).Value = 2621.4
$ws.Cells.Item(86, 13).Value = -939.8108000000002
$ws.Cells.Item(86, 14).Value = -4867.4

# BSM row 89
$ws.Cells.Item(89, 8).Value = 2129.3096
$ws.Cells.Item(89, 9).Value = 2062.8108
$ws.Cells.Item(89, 10).Value = 2621.4
$ws.Cells.Item(89, 11).Value = 10314.054
$ws.Cells.Item(89, 12).Value = 13107
$ws.Cells.Item(89, 13).Value = -4698.054
$ws.Cells.Item(89, 14).Value = -24339

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Cells.Item(31, 8).Value = 28392.51
$ws.Cells.Item(31, 10).Value = 4653.4
$ws.Cells.Item(31, 12).Value = 4653.4
$ws.Cells.Item(31, 14).Value = -5243.4

# CRP row 34
$ws.Cells.Item(34, 8).Value = 28392.51
$ws.Cells.Item(34, 10).Value = 4653.4
$ws.Cells.Item(34, 12).Value = 4653.4
$ws.Cells.Item(34, 14).Value = -5057.4

# CRP row 134
$ws.Cells.Item(134, 8).Value = 7685
$ws.Cells.Item(134, 9).Value = 5057.107
$ws.Cells.Item(134, 10).Value = 16882.625
$ws.Cells.Item(134, 11).Value = 15171.321
$ws.Cells.Item(134, 12).Value = 50647.875
$ws.Cells.Item(134, 13).Value = -12636.321
$ws.Cells.Item(134, 14).Value = -55717.875

$ws = $wb.Worksheets.Item("CUL")
# CUL row 5
$ws.Cells.Item(5, 8).Value = 2148.611
$ws.Cells.Item(5, 9).Value = 1366.6666
$ws.Cells.Item(5, 10).Value = 2219.697
$ws.Cells.Item(5, 11).Value = 4099.9998
$ws.Cells.Item(5, 12).Value = 6659.091
$ws.Cells.Item(5, 13).Value = -3987.9998
$ws.Cells.Item(5, 14).Value = -6883.091

# CUL row 56
$ws.Cells.Item(56, 8).Value = 4834.4443
$ws.Cells.Item(56, 9).Value = 4834.4443
$ws.Cells.Item(56, 11).Value = 4834.4443
$ws.Cells.Item(56, 13).Value = -4304.4443

# CUL row 112
$ws.Cells.Item(112, 8).Value = 3672.7273
$ws.Cells.Item(112, 10).Value = 4000
$ws.Cells.Item(112, 12).Value = 12000
$ws.Cells.Item(112, 14).Value = -14216

# CUL row 122
$ws.Cells.Item(122, 8).Value = 1814.2858
$ws.Cells.Item(122, 9).Value = 680
$ws.Cells.Item(122, 10).Value = 2444.4443
$ws.Cells.Item(122, 11).Value = 6120
$ws.Cells.Item(122, 12).Value = 21999.9987
$ws.Cells.Item(122, 13).Value = -3670
$ws.Cells.Item(122, 14).Value = -26899.9987

# CUL row 125
$ws.Cells.Item(125, 8).Value = 4732.5
$ws.Cells.Item(125, 10).Value = 7500
$ws.Cells.Item(125, 12).Value = 22500
$ws.Cells.Item(125, 14).Value = -32340

# CUL row 135
$ws.Cells.Item(135, 8).Value = 2148.611
$ws.Cells.Item(135, 9).Value = 1366.6666
$ws.Cells.Item(135, 10).Value = 2219.697
$ws.Cells.Item(135, 11).Value = 12299.9994
$ws.Cells.Item(135, 12).Value = 19977.273
$ws.Cells.Item(135, 13).Value = -9764.999400000001
$ws.Cells.Item(135, 14).Value = -25047.273

$ws = $wb.Worksheets.Item("GSM")
# GSM row 97
$ws.Cells.Item(97, 8).Value = 2207.6924
$ws.Cells.Item(97, 9).Value = 2461.25
$ws.Cells.Item(97, 10).Value = 1802
$ws.Cells.Item(97, 11).Value = 2461.25
$ws.Cells.Item(97, 12).Value = 1802
$ws.Cells.Item(97, 13).Value = -1965.25
$ws.Cells.Item(97, 14).Value = -2794

# GSM row 132
$ws.Cells.Item(132, 8).Value = 33911.2
$ws.Cells.Item(132, 9).Value = 40599.77
$ws.Cells.Item(132, 10).Value = 14588.667
$ws.Cells.Item(132, 11).Value = 121799.31
$ws.Cells.Item(132, 12).Value = 43766.001
$ws.Cells.Item(132, 13).Value = -119269.31
$ws.Cells.Item(132, 14).Value = -48826.001

$ws = $wb.Worksheets.Item("LTW")
# LTW row 46
$ws.Cells.Item(46, 8).Value = 1657.1428
$ws.Cells.Item(46, 9).Value = 1025
$ws.Cells.Item(46, 10).Value = 2500
$ws.Cells.Item(46, 11).Value = 1025
$ws.Cells.Item(46, 12).Value = 2500
$ws.Cells.Item(46, 13).Value = -837
$ws.Cells.Item(46, 14).Value = -2876

# LTW row 55
$ws.Cells.Item(55, 8).Value = 393.17142
$ws.Cells.Item(55, 9).Value = 334.14285
$ws.Cells.Item(55, 10).Value = 481.7143
$ws.Cells.Item(55, 11).Value = 334.14285
$ws.Cells.Item(55, 12).Value = 481.7143
$ws.Cells.Item(55, 13).Value = -161.14285
$ws.Cells.Item(55, 14).Value = -827.7143

# LTW row 93
$ws.Cells.Item(93, 8).Value = 1227.4445
$ws.Cells.Item(93, 9).Value = 1398.8889
$ws.Cells.Item(93, 10).Value = 1056
$ws.Cells.Item(93, 11).Value = 1398.8889
$ws.Cells.Item(93, 12).Value = 1056
$ws.Cells.Item(93, 13).Value = -150.8888999999999
$ws.Cells.Item(93, 14).Value = -3552

# LTW row 122
$ws.Cells.Item(122, 8).Value = 1877.3334
$ws.Cells.Item(122, 9).Value = 1907.1428
$ws.Cells.Item(122, 10).Value = 1851.25
$ws.Cells.Item(122, 11).Value = 5721.428400000001
$ws.Cells.Item(122, 12).Value = 5553.75
$ws.Cells.Item(122, 13).Value = -3271.428400000001
$ws.Cells.Item(122, 14).Value = -10453.75
